$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on every cell we touch so that numeric-looking
# strings (prices, and the hour value "14") are stored as text, matching
# the original inlineStr cell type used throughout this worksheet.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("D2").Value = "245.64"
$ws.Range("G2").Value = "14"

# Row 3
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "14"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("D4").Value = "5.307"
$ws.Range("G4").Value = "14"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05737"
$ws.Range("G5").Value = "14"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("D6").Value = "6.479"
$ws.Range("G6").Value = "14"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("D7").Value = "3.145"
$ws.Range("G7").Value = "14"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8169"
$ws.Range("G8").Value = "14"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8685"
$ws.Range("G9").Value = "14"

# Row 10
$ws.Range("B10").NumberFormat = "@"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1378"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("G10").Value = "14"

# Row 11
$ws.Range("B11").NumberFormat = "@"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.06986"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("G11").Value = "14"

# Row 12
$ws.Range("B12").NumberFormat = "@"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "0.03172"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G12").Value = "14"

# Row 13
$ws.Range("B13").NumberFormat = "@"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.02913"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("G13").Value = "14"

# Row 14
$ws.Range("B14").NumberFormat = "@"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09380"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("G14").Value = "14"

# Row 15
$ws.Range("B15").NumberFormat = "@"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "3.741"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("G15").Value = "14"

# Row 16
$ws.Range("B16").NumberFormat = "@"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "0.001534"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("G16").Value = "14"

# Row 17
$ws.Range("B17").NumberFormat = "@"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "0.04725"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("G17").Value = "14"

# Row 18
$ws.Range("B18").NumberFormat = "@"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "0.0006009"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("G18").Value = "14"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006180"
$ws.Range("G19").Value = "14"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001240"
$ws.Range("G20").Value = "14"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("D21").Value = "0.003861"
$ws.Range("G21").Value = "14"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00008798"
$ws.Range("G22").Value = "14"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("D23").Value = "3.537"
$ws.Range("G23").Value = "14"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("D24").Value = "2.139"
$ws.Range("G24").Value = "14"

# Row 25
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "14"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1328"
$ws.Range("G26").Value = "14"

# Row 27
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "14"

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0003014"
$ws.Range("E28").Value = "27UpBotsUBXTBestin24h"
$ws.Range("G28").Value = "14"

# Row 29
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "14"

# Row 30
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "14"

# Row 31
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "14"

# Row 32
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "14"

# Row 33
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "14"

# Row 34
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "14"

# Row 35
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "14"

# Row 36
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "14"

# Row 37
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "14"

# Row 38
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "14"

# Row 39
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "14"

# Row 40
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "14"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006439"
$ws.Range("G41").Value = "14"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1056"
$ws.Range("G42").Value = "14"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002200"
$ws.Range("G43").Value = "14"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007849"
$ws.Range("G44").Value = "14"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005205"
$ws.Range("G45").Value = "14"

# Row 46
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "14"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("D47").Value = "0.3499"
$ws.Range("G47").Value = "14"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("D48").Value = "0.001871"
$ws.Range("E48").Value = "47BOLOBOLO"
$ws.Range("G48").Value = "14"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("G49").Value = "14"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("G50").Value = "14"

# Row 51
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "14"
